$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New animation blocks appended below the existing "Gancho" data (rows 1-37
# already exist). We add four new animations: "Recibir daño" (1 frame),
# "Caerse" (3 frames), "Intentar levantarse" (2 frames, merged) and
# "Levantarse" (4 frames, merged) -- rows 38 through 47.
# ---------------------------------------------------------------------------

# --- Row 38: Recibir daño (single frame, bold + horizontal-center label) ---
$ws.Range("A38").Value = "Recibir daño"
$ws.Range("A38").Font.Bold = $true
$ws.Range("A38").HorizontalAlignment = -4108

$ws.Range("B38").Value = 1
$ws.Range("C38").Value = 26
$ws.Range("D38").Value = 712
$ws.Range("E38").Formula = "=80-C38"
$ws.Range("F38").Formula = "=795-D38"

# --- Rows 39-41: Caerse (3 frames, bold label with no alignment) ---
$ws.Range("A39").Value = "Caerse"
$ws.Range("A39").Font.Bold = $true

$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 112
$ws.Range("D39").Value = 700
$ws.Range("E39").Formula = "=172-C39"
$ws.Range("F39").Formula = "=795-D39"

$ws.Range("B40").Value = 2
$ws.Range("C40").Value = 214
$ws.Range("D40").Value = 704
$ws.Range("E40").Formula = "=308-C40"
$ws.Range("F40").Formula = "=756-D40"

$ws.Range("B41").Value = 3
$ws.Range("C41").Value = 337
$ws.Range("D41").Value = 753
$ws.Range("E41").Formula = "=453-C41"
$ws.Range("F41").Formula = "=789-D41"

# --- Rows 42-43: Intentar levantarse (2 frames, merged A42:A43) ---
$ws.Range("A42:A43").Merge()
$ws.Range("A42").Value = "Intentar levantarse"
$ws.Range("A42:A43").Font.Bold = $true
$ws.Range("A42:A43").HorizontalAlignment = -4108
$ws.Range("A42:A43").VerticalAlignment = -4108
$ws.Range("A42:A43").WrapText = $true

$ws.Range("B42").Value = 1
$ws.Range("C42").Value = 337
$ws.Range("D42").Value = 753
$ws.Range("E42").Formula = "=453-C42"
$ws.Range("F42").Formula = "=789-D42"

$ws.Range("B43").Value = 2
$ws.Range("C43").Value = 481
$ws.Range("D43").Value = 759
$ws.Range("E43").Formula = "=574-C43"
$ws.Range("F43").Formula = "=793-D43"

# --- Rows 44-47: Levantarse (4 frames, merged A44:A47) ---
$ws.Range("A44:A47").Merge()
$ws.Range("A44").Value = "Levantarse"
$ws.Range("A44:A47").Font.Bold = $true
$ws.Range("A44:A47").HorizontalAlignment = -4108
$ws.Range("A44:A47").VerticalAlignment = -4108

$ws.Range("B44").Value = 1
$ws.Range("C44").Value = 337
$ws.Range("D44").Value = 753
$ws.Range("E44").Formula = "=453-C44"
$ws.Range("F44").Formula = "=789-D44"

$ws.Range("B45").Value = 2
$ws.Range("C45").Value = 481
$ws.Range("D45").Value = 759
$ws.Range("E45").Formula = "=574-C45"
$ws.Range("F45").Formula = "=793-D45"

$ws.Range("B46").Value = 3
$ws.Range("C46").Value = 602
$ws.Range("D46").Value = 739
$ws.Range("E46").Formula = "=655-C46"
$ws.Range("F46").Formula = "=797-D46"

$ws.Range("B47").Value = 4
$ws.Range("C47").Value = 31
$ws.Range("D47").Value = 25
$ws.Range("E47").Value = 59
$ws.Range("F47").Value = 93

# ---------------------------------------------------------------------------
# Update the view: selection moves to the new last cell.
# ---------------------------------------------------------------------------
$ws.Range("F47").Select()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
